$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column D
$ws.Range("D2").Value = "canonical SMILES"

# Column D should get the same width as column C
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# For every data row (3 through 29), column D = column C value with
# bond-stereo slashes ("/" and "\") removed, producing the
# non-isomeric canonical SMILES.
for ($r = 3; $r -le 29; $r++) {
    $c = $ws.Cells.Item($r, 3).Value()
    if ($c -ne $null) {
        $d = $c.Replace("/", "").Replace("\", "")
        $ws.Cells.Item($r, 4).Value = $d
    }
}
